$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new log rows (8 and 9), matching the date style already used in B4:B7
$ws.Range("B7").Copy()
$ws.Range("B8:B9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B8").Value = (Get-Date -Year 2026 -Month 1 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D8").Value = "UI Start"

$ws.Range("B9").Value = (Get-Date -Year 2026 -Month 1 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D9").Value = "rules 0.5"
